$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "56.954.37"
$ws.Range("E2").Value = "  -0.54%  "

# Row 3
$ws.Range("D3").Value = "2.312.42"
$ws.Range("E3").Value = "  -1.58%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'530.64"
$ws.Range("E5").Value = "  +2.08%  "

# Row 6
$ws.Range("D6").Value = "'131.87"
$ws.Range("E6").Value = "  -3.29%  "

# Row 7
$ws.Range("E7").Value = "  -0.27%  "

# Row 8
$ws.Range("D8").Value = "'0.536"
$ws.Range("E8").Value = "  -0.31%  "

# Row 9
$ws.Range("D9").Value = "2.332.20"
$ws.Range("E9").Value = "  -1.20%  "

# Row 10
$ws.Range("D10").Value = "'0.101"
$ws.Range("E10").Value = "  -1.28%  "

# Row 11
$ws.Range("E11").Value = "  +0.24%  "

# Row 12
$ws.Range("D12").Value = "'5.28"
$ws.Range("E12").Value = "  -2.93%  "

# Row 13
$ws.Range("E13").Value = "  +0.26%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.730.40"
$ws.Range("E14").Value = "  -1.47%  "

# Row 15
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'23.42"
$ws.Range("E15").Value = "  -3.20%  "

# Row 16
$ws.Range("D16").Value = "56.999.64"
$ws.Range("E16").Value = "  -0.45%  "

# Row 17
$ws.Range("E17").Value = "  -2.13%  "

# Row 18
$ws.Range("D18").Value = "2.327.53"
$ws.Range("E18").Value = "  -1.55%  "

# Row 19
$ws.Range("D19").Value = "'337.95"
$ws.Range("E19").Value = "  +2.79%  "

# Row 20
$ws.Range("E20").Value = "  -1.51%  "

# Row 21
$ws.Range("D21").Value = "'6.91"
$ws.Range("E21").Value = "  +2.86%  "

# Row 22
$ws.Range("D22").Value = "'4.14"
$ws.Range("E22").Value = "  -2.15%  "

# Row 23
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.00%  "

# Row 24
$ws.Range("D24").Value = "'61.56"
$ws.Range("E24").Value = "  +0.61%  "

# Row 25
$ws.Range("E25").Value = "  +5.05%  "

# Row 26
$ws.Range("E26").Value = "  +0.34%  "

# Row 27
$ws.Range("D27").Value = "'0.996"
$ws.Range("E27").Value = "  -0.03%  "

# Row 28
$ws.Range("D28").Value = "'1.32"
$ws.Range("E28").Value = "  -0.68%  "

# Row 29
$ws.Range("D29").Value = "'170.51"
$ws.Range("E29").Value = "  -0.08%  "

# Row 30
$ws.Range("E30").Value = "  +0.86%  "

# Row 31
$ws.Range("D31").Value = "0.0₃0718"
$ws.Range("E31").Value = "  -3.26%  "

# Row 32
$ws.Range("D32").Value = "'6.08"
$ws.Range("E32").Value = "  -3.07%  "

# Row 33
$ws.Range("E33").Value = "  -0.41%  "

# Row 35
$ws.Range("D35").Value = "'0.992"
$ws.Range("E35").Value = "  -0.26%  "

# Row 36
$ws.Range("E36").Value = "  -2.82%  "

# Row 37
$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D37").Value = "'0.906"
$ws.Range("E37").Value = "  -1.73%  "

# Row 38
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'3.99"
$ws.Range("E38").Value = "  -0.96%  "

# Row 39
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "'39.03"
$ws.Range("E39").Value = "  +1.42%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'1.57"
$ws.Range("E40").Value = "  -0.14%  "

# Row 41
$ws.Range("D41").Value = "'148.28"
$ws.Range("E41").Value = "  -1.90%  "

# Row 42
$ws.Range("D42").Value = "'0.376"
$ws.Range("E42").Value = "  -1.70%  "

# Row 43
$ws.Range("E43").Value = "  -1.58%  "

# Row 44
$ws.Range("D44").Value = "'278.22"
$ws.Range("E44").Value = "  -1.33%  "

# Row 45
$ws.Range("D45").Value = "'5.05"
$ws.Range("E45").Value = "  -3.36%  "

# Row 46
$ws.Range("D46").Value = "'0.0927"
$ws.Range("E46").Value = "  -1.24%  "

# Row 47
$ws.Range("D47").Value = "'0.0502"
$ws.Range("E47").Value = "  -1.08%  "

# Row 48
$ws.Range("D48").Value = "'0.554"
$ws.Range("E48").Value = "  -1.46%  "

# Row 49
$ws.Range("D49").Value = "'18.49"
$ws.Range("E49").Value = "  +1.75%  "

# Row 50
$ws.Range("E50").Value = "  -2.10%  "

# Row 51
$ws.Range("E51").Value = "  -0.29%  "
